# Update existing match in row 2 with new odds data, and append a new
# match (row 7) for Suwon FC vs Seoul, South Korea K League 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: updated odds (columns G..BD correspond to columns 7..56) ----
$row2Values = @(
    1.91, 3.9, 3.7, 2.4, 2.5, 3.75, 1.02, 21, 1.14, 5.5,
    1.5, 2.63, 1.25, 3.75, 1.41, 2.62, 12, 12, 9, 19,
    13, 19, 21, 8, 11, 29, 81, 17, 23, 13,
    41, 23, 26, 4.5, 9.5, 15, 29, 41, 81, 3.75,
    7, 41, 251, 6, 17, 21, 51, 51, 101, 151
)

$col = 7
foreach ($v in $row2Values) {
    $ws.Cells.Item(2, $col).Value = $v
    $col = $col + 1
}

# ---- Row 7: brand-new match row ----
$ws.Cells.Item(7, 1).Value = "86Td3Gio"
$ws.Cells.Item(7, 2).Value = "26/10/2024"
$ws.Cells.Item(7, 3).Value = "04:30"
$ws.Cells.Item(7, 4).Value = "SOUTH KOREA - K LEAGUE 1"
$ws.Cells.Item(7, 5).Value = "Suwon FC"
$ws.Cells.Item(7, 6).Value = "Seoul"

$row7Values = @(
    3.1, 3.25, 2.3, 3.6, 2.2, 3, 1.05, 11, 1.29, 3.5,
    1.98, 1.88, 1.4, 2.75, 1.73, 2, 10, 15, 11, 34,
    23, 34, 10, 6.5, 13, 41, 201, 8.5, 11, 9.5,
    21, 19, 26, 5, 17, 26, 51, 67, 151, 2.75,
    7.5, 51, 501, 4.33, 13, 21, 41, 51, 151, 51
)

$col = 7
foreach ($v in $row7Values) {
    $ws.Cells.Item(7, $col).Value = $v
    $col = $col + 1
}
